$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43,8).Value = 2772.1072
$ws.Cells.Item(43,10).Value = 3374.125
$ws.Cells.Item(43,12).Value = 3374.125
$ws.Cells.Item(43,14).Value = -3512.125
$ws.Cells.Item(51,8).Value = 11999.667
$ws.Cells.Item(51,9).Value = 11999.667
$ws.Cells.Item(51,11).Value = 11999.667
$ws.Cells.Item(51,13).Value = -11515.667
$ws.Cells.Item(74,8).Value = 29466.166
$ws.Cells.Item(74,9).Value = 29466.166
$ws.Cells.Item(74,11).Value = 29466.166
$ws.Cells.Item(74,13).Value = -28530.166
$ws.Cells.Item(77,8).Value = 29466.166
$ws.Cells.Item(77,9).Value = 29466.166
$ws.Cells.Item(77,11).Value = 147330.83
$ws.Cells.Item(77,13).Value = -142650.83
$ws.Cells.Item(100,8).Value = 2694.3333
$ws.Cells.Item(100,9).Value = 1542
$ws.Cells.Item(100,11).Value = 1542
$ws.Cells.Item(100,13).Value = -1001
$ws.Cells.Item(135,8).Value = 1258.04
$ws.Cells.Item(135,9).Value = 915.9
$ws.Cells.Item(135,11).Value = 8243.1
$ws.Cells.Item(135,13).Value = -5708.1
$ws.Cells.Item(137,8).Value = 4009.652
$ws.Cells.Item(137,9).Value = 3522.3333
$ws.Cells.Item(137,10).Value = 4923.375
$ws.Cells.Item(137,11).Value = 10566.9999
$ws.Cells.Item(137,12).Value = 14770.125
$ws.Cells.Item(137,13).Value = -8016.999899999999
$ws.Cells.Item(137,14).Value = -19870.125
$ws.Cells.Item(138,8).Value = 4446.9653
$ws.Cells.Item(138,9).Value = 4819.091
$ws.Cells.Item(138,10).Value = 4321.0156
$ws.Cells.Item(138,11).Value = 14457.273
$ws.Cells.Item(138,12).Value = 12963.0468
$ws.Cells.Item(138,13).Value = -9317.273000000001
$ws.Cells.Item(138,14).Value = -23243.0468
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32,8).Value = 15324.361
$ws.Cells.Item(32,9).Value = 10739.953
$ws.Cells.Item(32,11).Value = 10739.953
$ws.Cells.Item(32,13).Value = -10452.953
$ws.Cells.Item(63,8).Value = 3397.111
$ws.Cells.Item(63,9).Value = 3321.75
$ws.Cells.Item(63,10).Value = 4000
$ws.Cells.Item(63,11).Value = 3321.75
$ws.Cells.Item(63,12).Value = 4000
$ws.Cells.Item(63,13).Value = -2635.75
$ws.Cells.Item(63,14).Value = -5372
$ws.Cells.Item(66,8).Value = 3397.111
$ws.Cells.Item(66,9).Value = 3321.75
$ws.Cells.Item(66,10).Value = 4000
$ws.Cells.Item(66,11).Value = 16608.75
$ws.Cells.Item(66,12).Value = 20000
$ws.Cells.Item(66,13).Value = -13176.75
$ws.Cells.Item(66,14).Value = -26864
$ws.Cells.Item(74,8).Value = 4566.0303
$ws.Cells.Item(74,9).Value = 2781.276
$ws.Cells.Item(74,11).Value = 2781.276
$ws.Cells.Item(74,13).Value = -1907.276
$ws.Cells.Item(77,8).Value = 4566.0303
$ws.Cells.Item(77,9).Value = 2781.276
$ws.Cells.Item(77,11).Value = 13906.38
$ws.Cells.Item(77,13).Value = -9538.379999999999
$ws.Cells.Item(120,8).Value = 66777
$ws.Cells.Item(120,10).Value = 66777
$ws.Cells.Item(120,12).Value = 66777
$ws.Cells.Item(120,14).Value = -76453
$ws.Cells.Item(125,8).Value = 96333.336
$ws.Cells.Item(125,10).Value = 96333.336
$ws.Cells.Item(125,12).Value = 96333.336
$ws.Cells.Item(125,14).Value = -106173.336
$ws.Cells.Item(128,8).Value = 0
$ws.Cells.Item(128,10).Value = 0
$ws.Cells.Item(128,12).Value = 0
$ws.Cells.Item(128,14).ClearContents()
$ws.Cells.Item(132,8).Value = 5666.8423
$ws.Cells.Item(132,9).Value = 2978
$ws.Cells.Item(132,10).Value = 15750
$ws.Cells.Item(132,11).Value = 8934
$ws.Cells.Item(132,12).Value = 47250
$ws.Cells.Item(132,13).Value = -6404
$ws.Cells.Item(132,14).Value = -52310
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94,8).Value = 2292.6875
$ws.Cells.Item(94,9).Value = 1390.4166
$ws.Cells.Item(94,11).Value = 1390.4166
$ws.Cells.Item(94,13).Value = -939.4166
$ws.Cells.Item(107,8).Value = 1677.9584
$ws.Cells.Item(107,9).Value = 1606.5
$ws.Cells.Item(107,11).Value = 1606.5
$ws.Cells.Item(107,13).Value = 313.5
$ws.Cells.Item(134,8).Value = 4443.4653
$ws.Cells.Item(134,9).Value = 3664.92
$ws.Cells.Item(134,10).Value = 9309.375
$ws.Cells.Item(134,11).Value = 10994.76
$ws.Cells.Item(134,12).Value = 27928.125
$ws.Cells.Item(134,13).Value = -8459.76
$ws.Cells.Item(134,14).Value = -32998.125
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31,8).Value = 7609.4
$ws.Cells.Item(31,9).Value = 3981.5862
$ws.Cells.Item(31,11).Value = 3981.5862
$ws.Cells.Item(31,13).Value = -3686.5862
$ws.Cells.Item(34,8).Value = 7609.4
$ws.Cells.Item(34,9).Value = 3981.5862
$ws.Cells.Item(34,11).Value = 3981.5862
$ws.Cells.Item(34,13).Value = -3779.5862
$ws.Cells.Item(62,8).Value = 3599.1667
$ws.Cells.Item(62,9).Value = 0
$ws.Cells.Item(62,11).Value = 0
$ws.Cells.Item(62,13).ClearContents()
$ws.Cells.Item(65,8).Value = 3599.1667
$ws.Cells.Item(65,9).Value = 0
$ws.Cells.Item(65,11).Value = 0
$ws.Cells.Item(65,13).ClearContents()
$ws.Cells.Item(93,8).Value = 8000
$ws.Cells.Item(93,9).Value = 8000
$ws.Cells.Item(93,11).Value = 8000
$ws.Cells.Item(93,13).Value = -6128
$ws.Cells.Item(122,8).Value = 4158.243
$ws.Cells.Item(122,10).Value = 4517.2144
$ws.Cells.Item(122,12).Value = 13551.6432
$ws.Cells.Item(122,14).Value = -18451.6432
$ws.Cells.Item(141,8).Value = 246117.1
$ws.Cells.Item(141,10).Value = 274120.4
$ws.Cells.Item(141,12).Value = 274120.4
$ws.Cells.Item(141,14).Value = -284480.4
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(7,8).Value = 463.6
$ws.Cells.Item(7,9).Value = 469
$ws.Cells.Item(7,11).Value = 1407
$ws.Cells.Item(7,13).Value = -1295
$ws.Cells.Item(74,8).Value = 5000
$ws.Cells.Item(74,10).Value = 5000
$ws.Cells.Item(74,12).Value = 15000
$ws.Cells.Item(74,14).Value = -17122
$ws.Cells.Item(77,8).Value = 5000
$ws.Cells.Item(77,10).Value = 5000
$ws.Cells.Item(77,12).Value = 45000
$ws.Cells.Item(77,14).Value = -55608
$ws.Cells.Item(107,8).Value = 493.57895
$ws.Cells.Item(107,9).Value = 326.7143
$ws.Cells.Item(107,11).Value = 980.1428999999999
$ws.Cells.Item(107,13).Value = 939.8571000000001
$ws.Cells.Item(128,8).Value = 254620.5
$ws.Cells.Item(128,9).Value = 254620.5
$ws.Cells.Item(128,11).Value = 763861.5
$ws.Cells.Item(128,13).Value = -758881.5
$ws.Cells.Item(132,8).Value = 3217.5
$ws.Cells.Item(132,9).Value = 3276
$ws.Cells.Item(132,10).Value = 3159
$ws.Cells.Item(132,11).Value = 29484
$ws.Cells.Item(132,12).Value = 28431
$ws.Cells.Item(132,13).Value = -26954
$ws.Cells.Item(132,14).Value = -33491
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80,8).Value = 1580
$ws.Cells.Item(80,9).Value = 1620
$ws.Cells.Item(80,10).Value = 1500
$ws.Cells.Item(80,11).Value = 1620
$ws.Cells.Item(80,12).Value = 1500
$ws.Cells.Item(80,13).Value = -622
$ws.Cells.Item(80,14).Value = -3496
$ws.Cells.Item(83,8).Value = 1580
$ws.Cells.Item(83,9).Value = 1620
$ws.Cells.Item(83,10).Value = 1500
$ws.Cells.Item(83,11).Value = 8100
$ws.Cells.Item(83,12).Value = 7500
$ws.Cells.Item(83,13).Value = -3108
$ws.Cells.Item(83,14).Value = -17484
$ws.Cells.Item(107,8).Value = 334.08334
$ws.Cells.Item(107,9).Value = 365.1111
$ws.Cells.Item(107,10).Value = 241
$ws.Cells.Item(107,11).Value = 365.1111
$ws.Cells.Item(107,12).Value = 241
$ws.Cells.Item(107,13).Value = 1554.8889
$ws.Cells.Item(107,14).Value = -4081
$ws.Cells.Item(138,8).Value = 203996.75
$ws.Cells.Item(138,10).Value = 203996.75
$ws.Cells.Item(138,12).Value = 203996.75
$ws.Cells.Item(138,14).Value = -214276.75
$ws.Cells.Item(140,8).Value = 0
$ws.Cells.Item(140,10).Value = 0
$ws.Cells.Item(140,12).Value = 0
$ws.Cells.Item(140,14).ClearContents()
$ws.Cells.Item(141,8).Value = 62118.31
$ws.Cells.Item(141,9).Value = 50389.668
$ws.Cells.Item(141,10).Value = 65636.89999999999
$ws.Cells.Item(141,11).Value = 50389.668
$ws.Cells.Item(141,12).Value = 65636.89999999999
$ws.Cells.Item(141,13).Value = -45209.668
$ws.Cells.Item(141,14).Value = -75996.89999999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(12,8).Value = 19155.8
$ws.Cells.Item(12,9).Value = 6000
$ws.Cells.Item(12,10).Value = 22444.75
$ws.Cells.Item(12,11).Value = 6000
$ws.Cells.Item(12,12).Value = 22444.75
$ws.Cells.Item(12,13).Value = -5830
$ws.Cells.Item(12,14).Value = -22784.75
$ws.Cells.Item(16,8).Value = 622
$ws.Cells.Item(16,10).Value = 852.4286
$ws.Cells.Item(16,12).Value = 852.4286
$ws.Cells.Item(16,14).Value = -1192.4286
$ws.Cells.Item(46,8).Value = 3813
$ws.Cells.Item(46,9).Value = 4525.5
$ws.Cells.Item(46,11).Value = 4525.5
$ws.Cells.Item(46,13).Value = -4337.5
$ws.Cells.Item(99,8).Value = 37410.625
$ws.Cells.Item(99,9).Value = 26666.666
$ws.Cells.Item(99,10).Value = 69642.5
$ws.Cells.Item(99,11).Value = 26666.666
$ws.Cells.Item(99,12).Value = 69642.5
$ws.Cells.Item(99,13).Value = -23671.666
$ws.Cells.Item(99,14).Value = -75632.5
$ws.Cells.Item(122,8).Value = 6344.3887
$ws.Cells.Item(122,9).Value = 5174.9165
$ws.Cells.Item(122,11).Value = 15524.7495
$ws.Cells.Item(122,13).Value = -13074.7495
$ws.Cells.Item(132,8).Value = 4656.125
$ws.Cells.Item(132,9).Value = 3519.9
$ws.Cells.Item(132,11).Value = 10559.7
$ws.Cells.Item(132,13).Value = -8029.700000000001
$ws.Cells.Item(136,8).Value = 6801.877
$ws.Cells.Item(136,9).Value = 3683.0435
$ws.Cells.Item(136,10).Value = 8911.677
$ws.Cells.Item(136,11).Value = 11049.1305
$ws.Cells.Item(136,12).Value = 26735.031
$ws.Cells.Item(136,13).Value = -8499.130500000001
$ws.Cells.Item(136,14).Value = -31835.031
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62,8).Value = 8599
$ws.Cells.Item(62,10).Value = 10998.333
$ws.Cells.Item(62,12).Value = 10998.333
$ws.Cells.Item(62,14).Value = -12246.333
$ws.Cells.Item(65,8).Value = 8599
$ws.Cells.Item(65,10).Value = 10998.333
$ws.Cells.Item(65,12).Value = 54991.665
$ws.Cells.Item(65,14).Value = -61231.665
$ws.Cells.Item(132,8).Value = 3941.3076
$ws.Cells.Item(132,9).Value = 2342.6155
$ws.Cells.Item(132,10).Value = 8737.385
$ws.Cells.Item(132,11).Value = 7027.8465
$ws.Cells.Item(132,12).Value = 26212.155
$ws.Cells.Item(132,13).Value = -4497.8465
$ws.Cells.Item(132,14).Value = -31272.155
$ws.Cells.Item(140,8).Value = 79665.664
$ws.Cells.Item(140,10).Value = 79500
$ws.Cells.Item(140,12).Value = 79500
$ws.Cells.Item(140,14).Value = -89860
$ws.Cells.Item(141,8).Value = 59999
$ws.Cells.Item(141,10).Value = 60000
$ws.Cells.Item(141,12).Value = 60000
$ws.Cells.Item(141,14).Value = -70360
